$wb = $excel.ActiveWorkbook

# --- Sheet "pro" (sheet1): update production values, column width, view ---
$wsPro = $wb.Worksheets.Item("pro")
$wsPro.Range("B2").Value = 5757204.42518376
$wsPro.Range("B3").Value = 5474195.5655403
$wsPro.Range("B4").Value = 5720925.031929665
$wsPro.Range("B5").Value = 6066031.490838885
$wsPro.Range("B6").Value = 5549876.12782229
$wsPro.Range("B7").Value = 5347076.173052053
$wsPro.Range("B8").Value = 5093605.721415728
$wsPro.Range("B9").Value = 4443302.16682553
$wsPro.Range("B10").Value = 4859429.617431719
$wsPro.Range("B11").Value = 4956218.006740789
$wsPro.Range("B12").Value = 4875524.650870575
$wsPro.Range("B13").Value = 4639010.530507527
$wsPro.Range("B14").Value = 4646467.162342135
$wsPro.Range("B15").Value = 5128224.741856335
$wsPro.Range("B16").Value = 5436748.16913763
$wsPro.Range("B17").Value = 5171681.046725798
$wsPro.Range("B18").Value = 6140549.787034335
$wsPro.Range("B19").Value = 6612581.09528717
$wsPro.Range("B20").Value = 7344950.480442911
$wsPro.Range("B21").Value = 7477710
$wsPro.Range("B22").Value = 7946350
$wsPro.Range("B23").Value = 8698557.856187724
$wsPro.Range("B24").Value = 9868150.234512335
$wsPro.Range("B25").Value = 10349411.1654146
$wsPro.Range("B26").Value = 10089763.786913153
$wsPro.Columns.Item(2).ColumnWidth = 14.8333

# --- Sheet "ind" (sheet2): update index values ---
$wsInd = $wb.Worksheets.Item("ind")
$wsInd.Range("B2").Value = 110.88270484554977
$wsInd.Range("B3").Value = 105.89059471901635
$wsInd.Range("B4").Value = 102.37862506126639
$wsInd.Range("B5").Value = 161.63470314676968
$wsInd.Range("B6").Value = 105.4320059411905
$wsInd.Range("B7").Value = 100.6852946731629
$wsInd.Range("B8").Value = 97.34596410455048
$wsInd.Range("B9").Value = 153.68917096863854
$wsInd.Range("B10").Value = 110.183970363138
$wsInd.Range("B11").Value = 105.22331833902248
$wsInd.Range("B12").Value = 101.73347958351262
$wsInd.Range("B13").Value = 160.61615168916993
$wsInd.Range("B14").Value = 106.9074948266841
$wsInd.Range("B15").Value = 111.26668320766107
$wsInd.Range("B16").Value = 104.25016470892207
$wsInd.Range("B17").Value = 114.95833468701943
$wsInd.Range("B18").Value = 105.3057667632521
$wsInd.Range("B19").Value = 100.7020231156873
$wsInd.Range("B20").Value = 95.31685066361952
$wsInd.Range("B21").Value = 101.22988836644618
$wsInd.Range("B22").Value = 106.94101798760872
$wsInd.Range("B23").Value = 113.82392962208158
$wsInd.Range("B24").Value = 115.52119164904492
$wsInd.Range("B25").Value = 123.3139505822282
$wsInd.Range("B26").Value = 110.59933179186326
$wsInd.Range("B27").Value = 105.93326024399906
$wsInd.Range("B28").Value = 104.2493379627183
$wsInd.Range("B29").Value = 101.67701579589966
$wsInd.Range("B30").Value = 97.8258560960867
$wsInd.Range("B31").Value = 98.31423800076786
$wsInd.Range("B32").Value = 95.78876351224146
$wsInd.Range("B33").Value = 107.22650310100215
$wsInd.Range("B34").Value = 101.67930083228798
$wsInd.Range("B35").Value = 99.25341596804628
$wsInd.Range("B36").Value = 97.54865945850273
$wsInd.Range("B37").Value = 98.26140157029855
$wsInd.Range("B38").Value = 95.68168719520851
$wsInd.Range("B39").Value = 95.16574695786984
$wsInd.Range("B40").Value = 98.81985367925783
$wsInd.Range("B41").Value = 104.61479490729513
$wsInd.Range("B42").Value = 99.5202201084785
$wsInd.Range("B43").Value = 105.31882371913512
$wsInd.Range("B44").Value = 104.27372655674934
$wsInd.Range("B45").Value = 110.00490454402623
$wsInd.Range("B46").Value = 112.72614522310162
$wsInd.Range("B47").Value = 108.7480480385009
$wsInd.Range("B48").Value = 106.16081546096096
$wsInd.Range("B49").Value = 119.34123508677727
$wsInd.Range("B50").Value = 120.24392566268186
$wsInd.Range("B51").Value = 124.56077256605113
$wsInd.Range("B52").Value = 130.93263807996968
$wsInd.Range("B53").Value = 138.34139509484086
$wsInd.Range("B54").Value = 120.75194701463552
$wsInd.Range("B55").Value = 122.47173987594401
$wsInd.Range("B56").Value = 128.3692224920335
$wsInd.Range("B57").Value = 142.38697618990386
$wsInd.Range("B58").Value = 132.93951830170624
$wsInd.Range("B59").Value = 138.6575814770154
$wsInd.Range("B60").Value = 137.39687635305214
$wsInd.Range("B61").Value = 135.5626492536738
$wsInd.Range("B62").Value = 112.52639095921172
$wsInd.Range("B63").Value = 104.77728334956889
$wsInd.Range("B64").Value = 134.64675933439818
$wsInd.Range("B65").Value = 154.51450010458035
$wsInd.Range("B66").Value = 162.50819799876845
$wsInd.Range("B67").Value = 172.20904919287423
$wsInd.Range("B68").Value = 168.16381968438665
$wsInd.Range("B69").Value = 180.1102970249527
$wsInd.Range("B70").Value = 181.08012734230567
$wsInd.Range("B71").Value = 189.82519905170818
$wsInd.Range("B72").Value = 191.76100689195155
$wsInd.Range("B73").Value = 209.2791716187585
$wsInd.Range("B74").Value = 204.3333586717344
$wsInd.Range("B75").Value = 218.73567656920747
$wsInd.Range("B76").Value = 222.82016453809038
$wsInd.Range("B77").Value = 219.65215490212609
$wsInd.Range("B78").Value = 216.06598817818102
$wsInd.Range("B79").Value = 232.38700023428802
$wsInd.Range("B80").Value = 231.27856009403462
$wsInd.Range("B81").Value = 236.01139676649876
$wsInd.Range("B82").Value = 226.19306279401144
$wsInd.Range("B83").Value = 225.5231197097944
$wsInd.Range("B84").Value = 221.97744221277432
$wsInd.Range("B85").Value = 231.53406874265573
$wsInd.Range("B86").Value = 284.2916475625974
$wsInd.Range("B87").Value = 273.2066041320614
$wsInd.Range("B88").Value = 273.0255982808133
$wsInd.Range("B89").Value = 277.50920993249946
$wsInd.Range("B90").Value = 284.3397960783906
$wsInd.Range("B91").Value = 253.66952260104742
$wsInd.Range("B92").Value = 269.99171907527733
$wsInd.Range("B93").Value = 268.5426450070523
$wsInd.Range("B94").Value = 318.0909535839056
$wsInd.Range("B95").Value = 266.83372862018734
$wsInd.Range("B96").Value = 298.30351576390376
$wsInd.Range("B97").Value = 289.3121129859738
$wsInd.Range("B98").Value = 350.09886977583034
$wsInd.Range("B99").Value = 277.36894830514603
$wsInd.Range("B100").Value = 324.15937671869085
$wsInd.Range("B101").Value = 308.4191930937124

# --- Sheet "conso" (sheet4): update consumption values ---
$wsConso = $wb.Worksheets.Item("conso")
$wsConso.Range("B2").Value = 2833561.4228866855
$wsConso.Range("B3").Value = 2694268.7426276295
$wsConso.Range("B4").Value = 2815706.2253714614
$wsConso.Range("B5").Value = 2985557.8699308033
$wsConso.Range("B6").Value = 2731517.8987471694
$wsConso.Range("B7").Value = 2631707.4747483395
$wsConso.Range("B8").Value = 2506952.215295386
$wsConso.Range("B9").Value = 2186889.2965822346
$wsConso.Range("B10").Value = 2391695.8620984363
$wsConso.Range("B11").Value = 2439333.1828502887
$wsConso.Range("B12").Value = 2399615.2856409606
$wsConso.Range("B13").Value = 2283210.8388029905
$wsConso.Range("B14").Value = 2286881.8910675156
$wsConso.Range("B15").Value = 2523993.2251696037
$wsConso.Range("B16").Value = 2675843.658777026
$wsConso.Range("B17").Value = 2545384.95723722
$wsConso.Range("B18").Value = 3022240.0553612798
$wsConso.Range("B19").Value = 3254562.7503222344
$wsConso.Range("B20").Value = 3615019.748091368
$wsConso.Range("B21").Value = 3680360
$wsConso.Range("B22").Value = 3619168
$wsConso.Range("B23").Value = 3833790.8603411024
$wsConso.Range("B24").Value = 4293845.894866913
$wsConso.Range("B25").Value = 4503252.949219136
$wsConso.Range("B26").Value = 4390274.751299919

# --- VA sheet (sheet3) recomputes automatically via formula =pro!Bx -conso!Bx ---

# --- Selection / active-sheet restore to match target view state ---
$wsInd.Activate()
$wsInd.Range("E19").Select()
$wsVA = $wb.Worksheets.Item("VA")
$wsVA.Activate()
$wsVA.Range("E19").Select()
$wsConso.Activate()
$wsConso.Range("E19").Select()
$wsPro.Activate()
$wsPro.Range("E19").Select()

